$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2163.4
$ws.Range("I40").Value = 1889
$ws.Range("J40").Value = 2575
$ws.Range("K40").Value = 1889
$ws.Range("L40").Value = 2575
$ws.Range("M40").Value = -1714
$ws.Range("N40").Value = -2925

# Row 53
$ws.Range("H53").Value = 870
$ws.Range("I53").Value = 929.61536
$ws.Range("J53").Value = 95
$ws.Range("K53").Value = 929.61536
$ws.Range("L53").Value = 95
$ws.Range("M53").Value = -292.61536
$ws.Range("N53").Value = -1369

# Row 64
$ws.Range("H64").Value = 3986.6667
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

# Row 67
$ws.Range("H67").Value = 3986.6667
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

# Row 132
$ws.Range("H132").Value = 8552648
$ws.Range("J132").Value = 3484.75
$ws.Range("L132").Value = 10454.25
$ws.Range("N132").Value = -15514.25

# Row 137
$ws.Range("H137").Value = 1786.2285
$ws.Range("I137").Value = 1625.8096
$ws.Range("K137").Value = 4877.4288
$ws.Range("M137").Value = -2327.4288

# Row 139
$ws.Range("H139").Value = 30930
$ws.Range("J139").Value = 34573.332
$ws.Range("L139").Value = 34573.332
$ws.Range("N139").Value = -44853.332

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 8099.4707
$ws.Range("I32").Value = 7484.357
$ws.Range("J32").Value = 9093.115
$ws.Range("K32").Value = 7484.357
$ws.Range("L32").Value = 9093.115
$ws.Range("M32").Value = -7197.357
$ws.Range("N32").Value = -9667.115

# Row 61
$ws.Range("H61").Value = 52632820
$ws.Range("I61").Value = 62501036
$ws.Range("K61").Value = 62501036
$ws.Range("M61").Value = -62500824

# Row 74
$ws.Range("H74").Value = 2975.182
$ws.Range("I74").Value = 2900
$ws.Range("J74").Value = 2982.7
$ws.Range("K74").Value = 2900
$ws.Range("L74").Value = 2982.7
$ws.Range("M74").Value = -2026
$ws.Range("N74").Value = -4730.7

# Row 77
$ws.Range("H77").Value = 2975.182
$ws.Range("I77").Value = 2900
$ws.Range("J77").Value = 2982.7
$ws.Range("K77").Value = 14500
$ws.Range("L77").Value = 14913.5
$ws.Range("M77").Value = -10132
$ws.Range("N77").Value = -23649.5

# Row 132
$ws.Range("H132").Value = 2438.578
$ws.Range("I132").Value = 2329.1904
$ws.Range("J132").Value = 2534.2917
$ws.Range("K132").Value = 6987.5712
$ws.Range("L132").Value = 7602.875100000001
$ws.Range("M132").Value = -4457.5712
$ws.Range("N132").Value = -12662.8751

# Row 136
$ws.Range("H136").Value = 52632820
$ws.Range("I136").Value = 62501036
$ws.Range("K136").Value = 187503108
$ws.Range("M136").Value = -187500558

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 3857.9
$ws.Range("I20").Value = 3739.8572
$ws.Range("J20").Value = 4133.3335
$ws.Range("K20").Value = 3739.8572
$ws.Range("L20").Value = 4133.3335
$ws.Range("M20").Value = -3492.8572
$ws.Range("N20").Value = -4627.3335

# Row 134
$ws.Range("H134").Value = 1788.6666
$ws.Range("I134").Value = 1454.8572
$ws.Range("K134").Value = 4364.571599999999
$ws.Range("M134").Value = -1829.571599999999

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1951.7576
$ws.Range("I31").Value = 1745.4193
$ws.Range("J31").Value = 5150
$ws.Range("K31").Value = 1745.4193
$ws.Range("L31").Value = 5150
$ws.Range("M31").Value = -1450.4193
$ws.Range("N31").Value = -5740

# Row 34
$ws.Range("H34").Value = 1951.7576
$ws.Range("I34").Value = 1745.4193
$ws.Range("J34").Value = 5150
$ws.Range("K34").Value = 1745.4193
$ws.Range("L34").Value = 5150
$ws.Range("M34").Value = -1543.4193
$ws.Range("N34").Value = -5554

# Row 58
$ws.Range("H58").Value = 3777.7046
$ws.Range("I58").Value = 1142.6666
$ws.Range("J58").Value = 7962.7646
$ws.Range("K58").Value = 1142.6666
$ws.Range("L58").Value = 7962.7646
$ws.Range("M58").Value = -939.6666
$ws.Range("N58").Value = -8368.7646

# Row 62
$ws.Range("H62").Value = 9526040
$ws.Range("I62").Value = 2346.875
$ws.Range("K62").Value = 2346.875
$ws.Range("M62").Value = -1722.875

# Row 65
$ws.Range("H65").Value = 9526040
$ws.Range("I65").Value = 2346.875
$ws.Range("K65").Value = 11734.375
$ws.Range("M65").Value = -8614.375

# Row 122
$ws.Range("H122").Value = 1111
$ws.Range("I122").Value = 1111
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3333
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -883
$ws.Range("N122").ClearContents()

# Row 135
$ws.Range("H135").Value = 33363.75
$ws.Range("J135").Value = 33363.75
$ws.Range("L135").Value = 33363.75
$ws.Range("N135").Value = -43503.75

# Row 136
$ws.Range("H136").Value = 3777.7046
$ws.Range("I136").Value = 1142.6666
$ws.Range("J136").Value = 7962.7646
$ws.Range("K136").Value = 3427.9998
$ws.Range("L136").Value = 23888.2938
$ws.Range("M136").Value = -877.9998000000001
$ws.Range("N136").Value = -28988.2938

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 680.4375
$ws.Range("I2").Value = 23
$ws.Range("J2").Value = 2126.8
$ws.Range("K2").Value = 138
$ws.Range("L2").Value = 12760.8
$ws.Range("M2").Value = -25
$ws.Range("N2").Value = -12986.8

# Row 5
$ws.Range("H5").Value = 714.7727
$ws.Range("I5").Value = 715
$ws.Range("K5").Value = 2145
$ws.Range("M5").Value = -2033

# Row 38
$ws.Range("H38").Value = 90.47059
$ws.Range("I38").Value = 61.384617
$ws.Range("J38").Value = 185
$ws.Range("K38").Value = 184.153851
$ws.Range("L38").Value = 555
$ws.Range("M38").Value = 162.846149
$ws.Range("N38").Value = -1249

# Row 131
$ws.Range("H131").Value = 15626391
$ws.Range("I131").Value = 142857680
$ws.Range("J131").Value = 1496.1404
$ws.Range("K131").Value = 428573040
$ws.Range("L131").Value = 4488.4212
$ws.Range("M131").Value = -428568000
$ws.Range("N131").Value = -14568.4212

# Row 135
$ws.Range("H135").Value = 714.7727
$ws.Range("I135").Value = 715
$ws.Range("K135").Value = 6435
$ws.Range("M135").Value = -3900

# Row 140
$ws.Range("H140").Value = 3049.6038
$ws.Range("I140").Value = 2550.3044
$ws.Range("J140").Value = 3432.4
$ws.Range("K140").Value = 7650.9132
$ws.Range("L140").Value = 10297.2
$ws.Range("M140").Value = -2470.9132
$ws.Range("N140").Value = -20657.2

# Row 141
$ws.Range("H141").Value = 58825600
$ws.Range("I141").Value = 66667852
$ws.Range("J141").Value = 8716.5
$ws.Range("K141").Value = 200003556
$ws.Range("L141").Value = 26149.5
$ws.Range("M141").Value = -199998376
$ws.Range("N141").Value = -36509.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 8383.333000000001
$ws.Range("I80").Value = 20000
$ws.Range("K80").Value = 20000
$ws.Range("M80").Value = -19002

# Row 83
$ws.Range("H83").Value = 8383.333000000001
$ws.Range("I83").Value = 20000
$ws.Range("K83").Value = 100000
$ws.Range("M83").Value = -95008

$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 1301.8334
$ws.Range("I68").Value = 1302
$ws.Range("J68").Value = 1301.6666
$ws.Range("K68").Value = 1302
$ws.Range("L68").Value = 1301.6666
$ws.Range("M68").Value = -553
$ws.Range("N68").Value = -2799.6666

# Row 71
$ws.Range("H71").Value = 1301.8334
$ws.Range("I71").Value = 1302
$ws.Range("J71").Value = 1301.6666
$ws.Range("K71").Value = 6510
$ws.Range("L71").Value = 6508.333000000001
$ws.Range("M71").Value = -2766
$ws.Range("N71").Value = -13996.333

# Row 100
$ws.Range("H100").Value = 988.2105
$ws.Range("I100").Value = 871.73334
$ws.Range("J100").Value = 1425
$ws.Range("K100").Value = 871.73334
$ws.Range("L100").Value = 1425
$ws.Range("M100").Value = -330.73334
$ws.Range("N100").Value = -2507

# Row 136
$ws.Range("H136").Value = 1940.7
$ws.Range("I136").Value = 1486.9286
$ws.Range("J136").Value = 2999.5
$ws.Range("K136").Value = 4460.7858
$ws.Range("L136").Value = 8998.5
$ws.Range("M136").Value = -1910.7858
$ws.Range("N136").Value = -14098.5

# Row 138
$ws.Range("H138").Value = 33698.57
$ws.Range("J138").Value = 33698.57
$ws.Range("L138").Value = 33698.57
$ws.Range("N138").Value = -43978.57

$ws = $wb.Worksheets.Item("WVR")
# Row 123
$ws.Range("H123").Value = 37500
$ws.Range("J123").Value = 37500
$ws.Range("L123").Value = 37500
$ws.Range("N123").Value = -47300
